$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Actualização do contacto do David João (row 3) with his phone number
$ws.Range("B3").Value = 919659339

# Reflect the cell last selected/clicked by the user when saving
$ws.Range("C14").Select()
